$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 226, shifting existing rows 226..269 down to 227..270.
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with the new weekly record.
$ws.Cells.Item(226, 1).Value = 7
$ws.Cells.Item(226, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(226, 3).Value = "Ñuble"
$ws.Cells.Item(226, 4).Value = 44694
$ws.Cells.Item(226, 5).Value = 16
$ws.Cells.Item(226, 6).Value = 100112023
$ws.Cells.Item(226, 7).Value = "Brócoli"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Segunda"
$ws.Cells.Item(226, 10).Value = 120
$ws.Cells.Item(226, 11).Value = 650
$ws.Cells.Item(226, 12).Value = 650
$ws.Cells.Item(226, 13).Value = 650
$ws.Cells.Item(226, 14).Value = "`$/unidad"
$ws.Cells.Item(226, 15).Value = "Región del Maule"
$ws.Cells.Item(226, 16).Value = 650
$ws.Cells.Item(226, 17).Value = 1
$ws.Cells.Item(226, 18).Value = "Hortaliza"
